# Update "want to go" counts (column F) on both the "展览" (exhibition)
# sheet and the "全部类型" (all types) sheet, which hold duplicated rows.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 7733
    $ws.Range("F3").Value = 7540
    $ws.Range("F12").Value = 104
    $ws.Range("F14").Value = 1118
}
